# Update to anaemia by intervention sheet
#
# The "Interventions anemia" worksheet had a blank spacer row (row 2)
# removed, which shifts all of the data (and the six existing cell
# comments) up by one row. In addition, a new "Outcome" header was added
# in A1 (matching the style already used by the rest of the header row),
# and the header row's height was set to a custom 14pt.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interventions anemia")

# --- 1. Delete the blank row 2; everything below shifts up by one row ---
$ws.Rows(2).Delete()

# --- 2. Add the new "Outcome" header cell in A1, matching the existing
#        header formatting (copy format from B1), and set the custom row
#        height for row 1 ---
$ws.Range("A1").Value = "Outcome"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows(1).RowHeight = 14

# --- 3. Re-anchor the six cell comments, which do not automatically move
#        when rows are deleted. Recreate each one row higher, preserving
#        its original text. ---
function Move-CommentTo($oldRef, $newRef, $text) {
    $old = $ws.Range($oldRef)
    if ($old.Comment -ne $null) {
        $old.Comment.Delete()
    }
    $ws.Range($newRef).AddComment($text)
}

Move-CommentTo "E16" "E15" "Ruth:`nmade this number up`n"
Move-CommentTo "G17" "G16" "Ruth:`nmader this number up"
Move-CommentTo "B18" "B17" " Janka Petravic:`nMutually exclusive with AMS, target = 1 - coverage of AMS - pregnant at risk of malaria not receiving IPTp or bednets"
Move-CommentTo "L18" "L17" "Ruth:`nmade this number up`n"
Move-CommentTo "B30" "B29" " Janka Petravic:`nAlso brestfeeding women up to 6 months?"
Move-CommentTo "B34" "B33" " Janka Petravic:`nThere is no unit cost for this to impact the budget."

# --- 4. Update the active selection to match the new state ---
$ws.Range("A3").Select()
